# This workbook's weekly "Ajo" (garlic) price log gets a new week's entry.
# A new row is inserted at row 286 (pushing the existing rows 286-334 down
# to 287-335), and the new row is populated with this week's figures. All
# other cells naturally shift down with the row insert, matching the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 286, shifting rows 286:334 -> 287:335
$ws.Rows.Item(286).Insert()

# Populate the newly inserted row 286 with the new week's record
$ws.Range("A286").Value = 7
$ws.Range("B286").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C286").Value = "Ñuble"
$ws.Range("D286").Value = 44951
$ws.Range("E286").Value = 16
$ws.Range("F286").Value = 100112003
$ws.Range("G286").Value = "Ajo"
$ws.Range("H286").Value = "Chino"
$ws.Range("I286").Value = "Primera"
$ws.Range("J286").Value = 60
$ws.Range("K286").Value = 16000
$ws.Range("L286").Value = 17000
$ws.Range("M286").Value = 16500
$ws.Range("N286").Value = "$/caja 10 kilos"
$ws.Range("O286").Value = "China"
$ws.Range("P286").Value = 1650
$ws.Range("Q286").Value = 10
$ws.Range("R286").Value = "Hortaliza"
